$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The hotel id "10002" is now prefixed to each existing 5-digit room/service
# id in column A (id_servicios), e.g. 10101 -> 1000210101.
$newIds = @(
    1000210101,
    1000210102,
    1000210103,
    1000210104,
    1000210105,
    1000210106,
    1000210107,
    1000210108,
    1000210109,
    1000210110
)

for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newIds[$i]
}

# Column A now holds longer numbers, so it was resized (best-fit) to
# comfortably show them.
$ws.Columns.Item(1).ColumnWidth = 10.33

# Leave the selection where it ended up after updating the id column.
$ws.Range("A2:A11").Select()
